# Handback report generation: refresh the "Latest HO Xliff Generate Date",
# "Correspond Handoff Datetime" and "Correspond Handback DateTime" timestamps
# that the CI handback-status report stamps onto each sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: G2 "Latest HO Xliff Generate Date" (also mirrored on de-de!H2)
$overview.Range("G2").Value = "2016-08-26 23:03:55"

# zh-cn sheet: H2 "Correspond Handoff Datetime", K2 "Correspond Handback DateTime"
$zhcn.Range("H2").Value = "2016-08-26 23:03:50"
$zhcn.Range("K2").Value = "2016-08-26 23:04:16"

# de-de sheet: H2 "Correspond Handoff Datetime" mirrors the Overview timestamp,
# K2 "Correspond Handback DateTime" gets its own refreshed stamp.
$dede.Range("H2").Value = "2016-08-26 23:03:55"
$dede.Range("K2").Value = "2016-08-26 23:04:23"
